$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws2 = $wb.Worksheets.Item("演出")
$ws3 = $wb.Worksheets.Item("本地生活")
$ws4 = $wb.Worksheets.Item("全部类型")

# --- 展览 (sheet1.xml) ---
$ws1.Range("F3").Value = 230
$ws1.Range("F4").Value = 483
$ws1.Range("F5").Value = 2069
$ws1.Range("F7").Value = 7849
$ws1.Range("F8").Value = 256
$ws1.Range("F9").Value = 37
$ws1.Range("F11").Value = 227
$ws1.Range("F12").Value = 1757
$ws1.Range("F13").Value = 1526
$ws1.Range("F15").Value = 166
$ws1.Range("F16").Value = 3903
$ws1.Range("F17").Value = 5956
$ws1.Range("F18").Value = 678
$ws1.Range("F20").Value = 1073
$ws1.Range("F21").Value = 1222
$ws1.Range("F22").Value = 411
$ws1.Range("F23").Value = 6156
$ws1.Range("F25").Value = 53
$ws1.Range("F26").Value = 4183
$ws1.Range("F28").Value = 1925
$ws1.Range("F29").Value = 1155
$ws1.Range("F30").Value = 293
$ws1.Range("F33").Value = 32
$ws1.Range("F35").Value = 37
$ws1.Range("F37").Value = 1145
$ws1.Range("F38").Value = 496
$ws1.Range("F39").Value = 1863
$ws1.Range("F40").Value = 96
$ws1.Range("F41").Value = 399
$ws1.Range("F42").Value = 148
$ws1.Range("F43").Value = 1123
$ws1.Range("F45").Value = 61
$ws1.Range("F48").Value = 166
$ws1.Range("F49").Value = 21

# --- 演出 (sheet2.xml) ---
$ws2.Range("F11").Value = 669
$ws2.Range("F12").Value = 360
$ws2.Range("F20").Value = 166
$ws2.Range("F22").Value = 73
$ws2.Range("G22").Value = 180
$ws2.Range("F27").Value = 125
$ws2.Range("F36").Value = 12

# --- 本地生活 (sheet3.xml) ---
$ws3.Range("F4").Value = 449
$ws3.Range("F7").Value = 465
$ws3.Range("F8").Value = 3081
$ws3.Range("F9").Value = 909
$ws3.Range("F10").Value = 1052
$ws3.Range("F11").Value = 1226
$ws3.Range("F12").Value = 1549

# --- 全部类型 (sheet4.xml) ---
$ws4.Range("F2").Value = 449
$ws4.Range("F4").Value = 230
$ws4.Range("F5").Value = 483
$ws4.Range("F6").Value = 465
$ws4.Range("F7").Value = 3081
$ws4.Range("F8").Value = 2069
$ws4.Range("F9").Value = 7851
$ws4.Range("F10").Value = 37
$ws4.Range("F11").Value = 909
$ws4.Range("F14").Value = 227
$ws4.Range("F15").Value = 1757
$ws4.Range("F16").Value = 1526
$ws4.Range("F17").Value = 1226
$ws4.Range("F19").Value = 669
$ws4.Range("F20").Value = 166
$ws4.Range("F21").Value = 1549
$ws4.Range("F22").Value = 3904
$ws4.Range("F23").Value = 360
$ws4.Range("F25").Value = 678
$ws4.Range("F27").Value = 1073
$ws4.Range("F28").Value = 1222
$ws4.Range("F29").Value = 411
$ws4.Range("F30").Value = 6156
$ws4.Range("F33").Value = 1925
$ws4.Range("F34").Value = 1155
$ws4.Range("F35").Value = 293
$ws4.Range("F36").Value = 32
$ws4.Range("F37").Value = 166
$ws4.Range("F39").Value = 73
$ws4.Range("G39").Value = 180
$ws4.Range("F40").Value = 496
$ws4.Range("F41").Value = 1863
$ws4.Range("F42").Value = 96
$ws4.Range("F43").Value = 399
$ws4.Range("F44").Value = 1123
$ws4.Range("F45").Value = 125
$ws4.Range("F49").Value = 166
